$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (values like "1.00", "611.47", "70.712.56" must not be coerced to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '70.712.56'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '3.523.93'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '611.47'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '173.70'
$ws.Range("E6").Value = '  +1.01%  '
$ws.Range("D7").Value = '3.517.76'
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("D11").Value = '7.39'
$ws.Range("E11").Value = '  +2.06%  '
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '46.55'
$ws.Range("E13").Value = '  -0.30%  '
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").Value = '4.095.55'
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '614.90'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = '3.523.63'
$ws.Range("E18").Value = '  -1.01%  '
$ws.Range("D19").Value = '70.758.10'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("D20").Value = '0.122'
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").Value = '17.80'
$ws.Range("E21").Value = '  +2.50%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").Value = '8.98'
$ws.Range("E23").Value = '  -4.95%  '
$ws.Range("D24").Value = '15.72'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '98.10'
$ws.Range("E25").Value = '  +1.31%  '
$ws.Range("D26").Value = '3.78'
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").Value = '2.60'
$ws.Range("E28").Value = '  -0.31%  '
$ws.Range("D29").Value = '33.76'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").Value = '9.14'
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").Value = '3.03'
$ws.Range("E31").Value = '  -1.18%  '
$ws.Range("E32").Value = '  -4.01%  '
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").Value = '615.60'
$ws.Range("E35").Value = '  +6.69%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("D37").Value = '10.85'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").Value = '3.52'
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").Value = '0.0475'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").Value = '57.01'
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = '3.375.27'
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").Value = '0.0₃0739'
$ws.Range("E44").Value = '  +4.81%  '
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("D48").Value = '2.57'
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("E49").Value = '  +0.22%  '
$ws.Range("D50").Value = '133.88'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("E51").Value = '  +0.00%  '
